$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(7, 9).Value = 'b'
$ws.Cells.Item(7, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(8, 9).Value = 'aa'
$ws.Cells.Item(8, 10).Value = 'Agree/Accept'
$ws.Cells.Item(34, 9).Value = 'sd'
$ws.Cells.Item(34, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(47, 9).Value = 'ba'
$ws.Cells.Item(47, 10).Value = 'Appreciation'
$ws.Cells.Item(62, 9).Value = 'sd'
$ws.Cells.Item(62, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(78, 9).Value = 'sv'
$ws.Cells.Item(78, 10).Value = 'Statement-opinion'
$ws.Cells.Item(81, 9).Value = 'sv'
$ws.Cells.Item(81, 10).Value = 'Statement-opinion'
$ws.Cells.Item(89, 9).Value = 'b'
$ws.Cells.Item(89, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(93, 9).Value = 'sd'
$ws.Cells.Item(93, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(100, 9).Value = 'ba'
$ws.Cells.Item(100, 10).Value = 'Appreciation'
$ws.Cells.Item(117, 9).Value = 'b'
$ws.Cells.Item(117, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(139, 9).Value = 'ba'
$ws.Cells.Item(139, 10).Value = 'Appreciation'
$ws.Cells.Item(158, 9).Value = 'ba'
$ws.Cells.Item(158, 10).Value = 'Appreciation'
$ws.Cells.Item(226, 9).Value = 'b'
$ws.Cells.Item(226, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(250, 9).Value = 'ba'
$ws.Cells.Item(250, 10).Value = 'Appreciation'
$ws.Cells.Item(273, 9).Value = 'sd'
$ws.Cells.Item(273, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(279, 9).Value = 'aa'
$ws.Cells.Item(279, 10).Value = 'Agree/Accept'
$ws.Cells.Item(281, 9).Value = 'sd'
$ws.Cells.Item(281, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(305, 9).Value = 'aa'
$ws.Cells.Item(305, 10).Value = 'Agree/Accept'
$ws.Cells.Item(311, 9).Value = 'b'
$ws.Cells.Item(311, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(317, 9).Value = 'sv'
$ws.Cells.Item(317, 10).Value = 'Statement-opinion'
$ws.Cells.Item(335, 9).Value = 'b'
$ws.Cells.Item(335, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(337, 9).Value = 'sv'
$ws.Cells.Item(337, 10).Value = 'Statement-opinion'
$ws.Cells.Item(339, 9).Value = 'sd'
$ws.Cells.Item(339, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(340, 9).Value = 'aa'
$ws.Cells.Item(340, 10).Value = 'Agree/Accept'
$ws.Cells.Item(344, 9).Value = 'aa'
$ws.Cells.Item(344, 10).Value = 'Agree/Accept'
$ws.Cells.Item(345, 9).Value = 'aa'
$ws.Cells.Item(345, 10).Value = 'Agree/Accept'
$ws.Cells.Item(349, 9).Value = 'qy'
$ws.Cells.Item(349, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(351, 9).Value = 'aa'
$ws.Cells.Item(351, 10).Value = 'Agree/Accept'
$ws.Cells.Item(352, 9).Value = 'aa'
$ws.Cells.Item(352, 10).Value = 'Agree/Accept'
$ws.Cells.Item(353, 9).Value = 'b'
$ws.Cells.Item(353, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(354, 9).Value = 'sd'
$ws.Cells.Item(354, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(357, 9).Value = 'b'
$ws.Cells.Item(357, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(359, 9).Value = 'sd'
$ws.Cells.Item(359, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(360, 9).Value = 'ba'
$ws.Cells.Item(360, 10).Value = 'Appreciation'
$ws.Cells.Item(361, 9).Value = 'aa'
$ws.Cells.Item(361, 10).Value = 'Agree/Accept'
$ws.Cells.Item(370, 9).Value = 'ba'
$ws.Cells.Item(370, 10).Value = 'Appreciation'
$ws.Cells.Item(374, 9).Value = 'sd'
$ws.Cells.Item(374, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(375, 9).Value = 'sd'
$ws.Cells.Item(375, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(376, 9).Value = 'ba'
$ws.Cells.Item(376, 10).Value = 'Appreciation'
$ws.Cells.Item(385, 9).Value = 'ba'
$ws.Cells.Item(385, 10).Value = 'Appreciation'
$ws.Cells.Item(391, 9).Value = 'aa'
$ws.Cells.Item(391, 10).Value = 'Agree/Accept'
$ws.Cells.Item(392, 9).Value = 'aa'
$ws.Cells.Item(392, 10).Value = 'Agree/Accept'
$ws.Cells.Item(398, 9).Value = 'b'
$ws.Cells.Item(398, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(413, 9).Value = 'ba'
$ws.Cells.Item(413, 10).Value = 'Appreciation'
$ws.Cells.Item(417, 9).Value = 'sd'
$ws.Cells.Item(417, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(422, 9).Value = 'b'
$ws.Cells.Item(422, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(424, 9).Value = 'sd'
$ws.Cells.Item(424, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(433, 9).Value = 'sv'
$ws.Cells.Item(433, 10).Value = 'Statement-opinion'
$ws.Cells.Item(434, 9).Value = 'sd'
$ws.Cells.Item(434, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(435, 9).Value = 'sd'
$ws.Cells.Item(435, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(440, 9).Value = '%'
$ws.Cells.Item(440, 10).Value = 'Uninterpretable'
$ws.Cells.Item(453, 9).Value = 'sd'
$ws.Cells.Item(453, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(484, 9).Value = 'aa'
$ws.Cells.Item(484, 10).Value = 'Agree/Accept'
$ws.Cells.Item(488, 9).Value = 'sv'
$ws.Cells.Item(488, 10).Value = 'Statement-opinion'
$ws.Cells.Item(490, 9).Value = 'sv'
$ws.Cells.Item(490, 10).Value = 'Statement-opinion'
$ws.Cells.Item(552, 9).Value = 'sd'
$ws.Cells.Item(552, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(559, 9).Value = '%'
$ws.Cells.Item(559, 10).Value = 'Uninterpretable'
$ws.Cells.Item(573, 9).Value = 'sd'
$ws.Cells.Item(573, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(580, 9).Value = 'sv'
$ws.Cells.Item(580, 10).Value = 'Statement-opinion'
$ws.Cells.Item(581, 9).Value = 'b'
$ws.Cells.Item(581, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(585, 9).Value = 'b'
$ws.Cells.Item(585, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(598, 9).Value = 'sd'
$ws.Cells.Item(598, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(606, 9).Value = 'sd'
$ws.Cells.Item(606, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(609, 9).Value = 'b'
$ws.Cells.Item(609, 10).Value = 'Acknowledge (Backchannel)'
